# Updated cryptos list on Thu Jan 18 08:12:54 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.913.32"
$ws.Range("E2").Value = "  +0.79%  "

$ws.Range("D3").Value = "2.536.38"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "311.84"
$ws.Range("E5").Value = "  +1.00%  "

$ws.Range("D6").Value = "100.88"
$ws.Range("E6").Value = "  +3.85%  "

$ws.Range("E7").Value = "  -0.54%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").Value = "0.525"
$ws.Range("E9").Value = "  -0.36%  "

$ws.Range("D10").Value = "35.84"
$ws.Range("E10").Value = "  +1.11%  "

$ws.Range("D11").Value = "0.0807"
$ws.Range("E11").Value = "  +0.44%  "

$ws.Range("D12").Value = "7.33"
$ws.Range("E12").Value = "  -0.40%  "

$ws.Range("E13").Value = "  +2.13%  "

$ws.Range("D14").Value = "2.926.94"
$ws.Range("E14").Value = "  -0.61%  "

$ws.Range("D15").Value = "15.41"
$ws.Range("E15").Value = "  -1.42%  "

$ws.Range("D16").Value = "2.544.94"
$ws.Range("E16").Value = "  -1.05%  "

$ws.Range("D17").Value = "0.819"
$ws.Range("E17").Value = "  -1.40%  "

$ws.Range("D18").Value = "42.891.59"
$ws.Range("E18").Value = "  +0.67%  "

$ws.Range("D19").Value = "6.69"
$ws.Range("E19").Value = "  -0.13%  "

$ws.Range("D20").Value = "12.36"
$ws.Range("E20").Value = "  +0.35%  "

$subscriptThree = [char]0x2083
$ws.Range("D21").Value = "0.0$($subscriptThree)0955"
$ws.Range("E21").Value = "  +0.20%  "

$ws.Range("D22").Value = "69.99"
$ws.Range("E22").Value = "  +1.24%  "

$ws.Range("D23").Value = "244.01"
$ws.Range("E23").Value = "  -0.76%  "

$ws.Range("D24").Value = "2.89"
$ws.Range("E24").Value = "  -0.39%  "

$ws.Range("E25").Value = "  +0.62%  "

$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").Value = "25.58"
$ws.Range("E27").Value = "  -3.14%  "

$ws.Range("E28").Value = "  -1.33%  "

$ws.Range("D29").Value = "10.22"
$ws.Range("E29").Value = "  +1.26%  "

$ws.Range("D30").Value = "38.70"
$ws.Range("E30").Value = "  -2.46%  "

$ws.Range("D31").Value = "159.08"
$ws.Range("E31").Value = "  +0.88%  "

$ws.Range("E32").Value = "  +2.88%  "

$ws.Range("D33").Value = "2.80"
$ws.Range("E33").Value = "  +8.38%  "

$ws.Range("E34").Value = "  +2.71%  "

$ws.Range("D35").Value = "0.0794"
$ws.Range("E35").Value = "  +0.61%  "

$ws.Range("D36").Value = "18.38"
$ws.Range("E36").Value = "  +0.02%  "

$ws.Range("E37").Value = "  -3.94%  "

$ws.Range("E38").Value = "  -4.11%  "

$ws.Range("E39").Value = "  +0.93%  "

$ws.Range("E40").Value = "  +0.43%  "

$ws.Range("E41").Value = "  +4.24%  "

$ws.Range("D42").Value = "21.99"
$ws.Range("E42").Value = "  -2.15%  "

# Rows 43 and 44 swap: FirstDigitalUSD moves up to rank 41, NEARProtocol moves to rank 42
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.19%  "

$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "3.33"
$ws.Range("E44").Value = "  +4.86%  "

$ws.Range("E45").Value = "  +0.50%  "

$ws.Range("D46").Value = "2.001.70"
$ws.Range("E46").Value = "  +0.98%  "

$ws.Range("D47").Value = "9.13"
$ws.Range("E47").Value = "  +2.19%  "

$ws.Range("D48").Value = "2.779.02"
$ws.Range("E48").Value = "  -0.64%  "

$ws.Range("E49").Value = "  +0.52%  "

$ws.Range("D50").Value = "80.16"
$ws.Range("E50").Value = "  -0.17%  "

$ws.Range("D51").Value = "72.56"
$ws.Range("E51").Value = "  -0.53%  "
